{"js": "// Remove the extraneous space before the comma in\n// \"... maiores matem\u00e1ticos , foi ...\" -> \"... maiores matem\u00e1ticos, foi ...\"\n// (third line of the bullet paragraph \u2014 see commit message\n// \"Removi um espa\u00e7o na terceira linha\").\nconst body = context.document.body;\n\nconst results = body.search(\"matem\u00e1ticos , foi\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found in document body.\");\n}\n\n// Replace the matched text (which still contains the stray space before\n// the comma) with the corrected text that has the space removed.\nresults.items[0].insertText(\"matem\u00e1ticos, foi\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Remove the extraneous space before the comma in\n# \"... maiores matem\u00e1ticos , foi ...\" -> \"... maiores matem\u00e1ticos, foi ...\"\n# (third line of the bullet paragraph \u2014 see commit message\n# \"Removi um espa\u00e7o na terceira linha\").\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"matem\u00e1ticos , foi\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"matem\u00e1ticos, foi\"\n\n$found = $find.Execute(\n    [ref]$find.Text,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]$find.Replacement.Text,\n    [ref]2\n)\n\nif (-not $found) {\n    throw \"Target phrase 'matem\u00e1ticos , foi' was not found in the document.\"\n}\n"}
